$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConsignmentID")

# Existing rows 16 & 17: spinner value "Y" -> "N"
$ws.Range("A16").Value = "N"
$ws.Range("A17").Value = "N"

# New rows 18-23: additional ConsignmentDetails test data ("Handled spinner and Billet To")
# Row 23 is populated before row 22 to mirror the original authoring order
# (reflected in the shared-string table insertion order).
$rows = @(
    @{ Row = 18; A = "N"; B = "UAT42092138"; E = "S"; F = "Surekha"; G = 123; H = "aaaaaa"; I = "qqqqqqq" },
    @{ Row = 19; A = "N"; B = "UAT42092140"; E = "A"; F = "Anand";   G = 456; H = "wwww";   I = "eeeeeee" },
    @{ Row = 20; A = "N"; B = "UAT42092141"; E = "S"; F = "Surekha"; G = 123; H = "aaaaaa"; I = "qqqqqqq" },
    @{ Row = 21; A = "Y"; B = "UAT42092142"; E = "A"; F = "Anand";   G = 456; H = "wwww";   I = "eeeeeee" },
    @{ Row = 23; A = "N"; B = "UAT42092144"; E = "A"; F = "Anand";   G = 456; H = "wwww";   I = "eeeeeee" },
    @{ Row = 22; A = "Y"; B = "UAT42092143"; E = "S"; F = "Surekha"; G = 123; H = "aaaaaa"; I = "qqqqqqq" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = "CELW01"
    $ws.Cells.Item($row, 4).Value = "Cell"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = 1
    $ws.Cells.Item($row, 11).Value = "SP12345678"
}

$ws.Range("A22").Select()
